$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 579.8
$ws.Range("J39").Value = 2499.5
$ws.Range("L39").Value = 7498.5
$ws.Range("N39").Value = -8090.5

$ws.Range("H51").Value = 3989.1738
$ws.Range("I51").Value = 4388.375
$ws.Range("K51").Value = 4388.375
$ws.Range("M51").Value = -3904.375

$ws.Range("H94").Value = 742.5
$ws.Range("I94").Value = 742.5
$ws.Range("K94").Value = 742.5
$ws.Range("M94").Value = -291.5

$ws.Range("H99").Value = 809.1429000000001
$ws.Range("I99").Value = 1582
$ws.Range("J99").Value = 500
$ws.Range("K99").Value = 4746
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = -3248
$ws.Range("N99").Value = -4496

$ws.Range("H101").Value = 672.125
$ws.Range("I101").Value = 479.5
$ws.Range("K101").Value = 1438.5
$ws.Range("M101").Value = 183.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3061
$ws.Range("I45").Value = 2052.1333
$ws.Range("J45").Value = 3857.4736
$ws.Range("K45").Value = 2052.1333
$ws.Range("L45").Value = 3857.4736
$ws.Range("M45").Value = -1675.1333
$ws.Range("N45").Value = -4611.473599999999

$ws.Range("H88").Value = 1394.125
$ws.Range("I88").Value = 1797.5
$ws.Range("K88").Value = 1797.5
$ws.Range("M88").Value = -1391.5

$ws.Range("H91").Value = 1394.125
$ws.Range("I91").Value = 1797.5
$ws.Range("K91").Value = 1797.5
$ws.Range("M91").Value = -393.5

$ws.Range("H102").Value = 1815.3529
$ws.Range("I102").Value = 1557.4
$ws.Range("K102").Value = 1557.4
$ws.Range("M102").Value = 64.59999999999991

$ws.Range("H110").Value = 1928.32
$ws.Range("I110").Value = 1691.3182
$ws.Range("K110").Value = 1691.3182
$ws.Range("M110").Value = 353.6818000000001

$ws.Range("H122").Value = 6035.8486
$ws.Range("I122").Value = 4952.2666
$ws.Range("K122").Value = 14856.7998
$ws.Range("M122").Value = -12406.7998

$ws.Range("H132").Value = 3019.72
$ws.Range("I132").Value = 2978.875
$ws.Range("K132").Value = 8936.625
$ws.Range("M132").Value = -6406.625

$ws.Range("H133").Value = 199000
$ws.Range("J133").Value = 199000
$ws.Range("L133").Value = 199000
$ws.Range("N133").Value = -204060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 285020
$ws.Range("J42").Value = 285020
$ws.Range("L42").Value = 285020
$ws.Range("N42").Value = -285676

$ws.Range("H134").Value = 3379.389
$ws.Range("I134").Value = 3239.9375
$ws.Range("K134").Value = 9719.8125
$ws.Range("M134").Value = -7184.8125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 505
$ws.Range("I22").Value = 505
$ws.Range("K22").Value = 505
$ws.Range("M22").Value = -155

$ws.Range("H99").Value = 3356.4285
$ws.Range("I99").Value = 2999
$ws.Range("K99").Value = 2999
$ws.Range("M99").Value = -1501

$ws.Range("H114").Value = 49778
$ws.Range("J114").Value = 49778
$ws.Range("L114").Value = 49778
$ws.Range("N114").Value = -58456

$ws.Range("H122").Value = 4946.4707
$ws.Range("I122").Value = 4343.5713
$ws.Range("J122").Value = 5368.5
$ws.Range("K122").Value = 13030.7139
$ws.Range("L122").Value = 16105.5
$ws.Range("M122").Value = -10580.7139
$ws.Range("N122").Value = -21005.5

$ws.Range("H126").Value = 3356.4285
$ws.Range("I126").Value = 2999
$ws.Range("K126").Value = 8997
$ws.Range("M126").Value = -6527

$ws.Range("H132").Value = 3992.8125
$ws.Range("I132").Value = 2891.4443
$ws.Range("K132").Value = 8674.332900000001
$ws.Range("M132").Value = -6144.332900000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 9529.583000000001
$ws.Range("I14").Value = 9529.583000000001
$ws.Range("K14").Value = 28588.749
$ws.Range("M14").Value = -28415.749

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 4039998.5
$ws.Range("I24").Value = 20000000
$ws.Range("K24").Value = 20000000
$ws.Range("M24").Value = -19999827

$ws.Range("H80").Value = 4000
$ws.Range("I80").Value = 4000
$ws.Range("K80").Value = 4000
$ws.Range("M80").Value = -3002

$ws.Range("H83").Value = 4000
$ws.Range("I83").Value = 4000
$ws.Range("K83").Value = 20000
$ws.Range("M83").Value = -15008

$ws.Range("H121").Value = 86662.664
$ws.Range("J121").Value = 86662.664
$ws.Range("L121").Value = 86662.664
$ws.Range("N121").Value = -90156.664

$ws.Range("H122").Value = 2059.4
$ws.Range("I122").Value = 1991.4615
$ws.Range("J122").Value = 2185.5715
$ws.Range("K122").Value = 5974.3845
$ws.Range("L122").Value = 6556.7145
$ws.Range("M122").Value = -3524.3845
$ws.Range("N122").Value = -11456.7145

$ws.Range("H126").Value = 3461.1538

$ws.Range("H132").Value = 3035.0264
$ws.Range("I132").Value = 2640.0417
$ws.Range("K132").Value = 7920.125100000001
$ws.Range("M132").Value = -5390.125100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6020.2085
$ws.Range("I7").Value = 6626.364
$ws.Range("K7").Value = 6626.364
$ws.Range("M7").Value = -6514.364

$ws.Range("H40").Value = 2386.25
$ws.Range("I40").Value = 2301.5715
$ws.Range("K40").Value = 2301.5715
$ws.Range("M40").Value = -2165.5715

$ws.Range("H82").Value = 3023.611
$ws.Range("J82").Value = 3103.8
$ws.Range("L82").Value = 3103.8
$ws.Range("N82").Value = -3825.8

$ws.Range("H85").Value = 3023.611
$ws.Range("J85").Value = 3103.8
$ws.Range("L85").Value = 3103.8
$ws.Range("N85").Value = -5599.8

$ws.Range("H93").Value = 100002344
$ws.Range("I93").Value = 200001620
$ws.Range("J93").Value = 3059.8
$ws.Range("K93").Value = 200001620
$ws.Range("L93").Value = 3059.8
$ws.Range("M93").Value = -200000372
$ws.Range("N93").Value = -5555.8

$ws.Range("H122").Value = 9847.686
$ws.Range("I122").Value = 9011.727999999999
$ws.Range("J122").Value = 11262.385
$ws.Range("K122").Value = 27035.184
$ws.Range("L122").Value = 33787.155
$ws.Range("M122").Value = -24585.184
$ws.Range("N122").Value = -38687.155

$ws.Range("H126").Value = 6020.2085
$ws.Range("I126").Value = 6626.364
$ws.Range("K126").Value = 19879.092
$ws.Range("M126").Value = -17409.092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4431.7856
$ws.Range("I81").Value = 3464
$ws.Range("J81").Value = 4818.9
$ws.Range("K81").Value = 6928
$ws.Range("L81").Value = 9637.799999999999
$ws.Range("M81").Value = -5867
$ws.Range("N81").Value = -11759.8

$ws.Range("H84").Value = 4431.7856
$ws.Range("I84").Value = 3464
$ws.Range("J84").Value = 4818.9
$ws.Range("K84").Value = 34640
$ws.Range("L84").Value = 48189
$ws.Range("M84").Value = -29336
$ws.Range("N84").Value = -58797

$ws.Range("H122").Value = 9913.1
$ws.Range("I122").Value = 10506.143
$ws.Range("K122").Value = 31518.429
$ws.Range("M122").Value = -29068.429

$ws.Range("H126").Value = 3599.3333
$ws.Range("I126").Value = 2899.25
$ws.Range("K126").Value = 8697.75
$ws.Range("M126").Value = -6227.75
